# Atualiza os dados da BIBI (faturamento_anual) para o ano de 2025 (linha 9)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value = 3231969.27
$ws.Range("C9").Value = 501581.99
$ws.Range("D9").Value = 3733551.26
$ws.Range("E9").Value = 13.43444766310775
$ws.Range("F9").Value = 86.56555233689225
$ws.Range("G9").Value = -51.52460166960253
$ws.Range("H9").Value = -41.63510233975794
$ws.Range("I9").Value = 32196
$ws.Range("J9").Value = 1359
$ws.Range("K9").Value = 33555
$ws.Range("L9").Value = 23164
$ws.Range("M9").Value = 161.1790390260749
$ws.Range("N9").Value = 10.04021279415479
